# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> used by the (single) Slide Master - currently "Integral"
#   ppt/theme/theme2.xml  -> used by the Notes Master           - currently "Office Theme"
# The authored edit swaps their colour schemes: theme1 becomes the stock
# "Office" palette, theme2 becomes "Integral". The color scheme is reached
# through Slide.ThemeColorScheme (it edits the theme backing the slide
# master, i.e. theme1.xml) - PowerPoint's RGB() is 0x00BBGGRR (OLE_COLOR),
# so convert each target hex RRGGBB accordingly.

function ToOleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation

# Target palette = the stock Office theme colours, in
# dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink order.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = ToOleColor $officeColors[$i - 1]
}
